# Updates the cryptos price/volume table to the latest scraped values
# (refresh performed by the scheduled GitHub Actions scraper).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.301.62"
$ws.Range("E2").Value = "'  -0.25%  "
$ws.Range("D3").Value = "'3.138.29"
$ws.Range("E3").Value = "'  -1.23%  "
$ws.Range("E4").Value = "'  -0.05%  "
$ws.Range("D5").Value = "'572.36"
$ws.Range("E5").Value = "'  +0.21%  "
$ws.Range("D6").Value = "'163.79"
$ws.Range("E6").Value = "'  -3.07%  "
$ws.Range("E7").Value = "'  -0.09%  "
$ws.Range("D8").Value = "'0.575"
$ws.Range("E8").Value = "'  -5.21%  "
$ws.Range("D9").Value = "'3.147.65"
$ws.Range("E9").Value = "'  -1.21%  "
$ws.Range("E10").Value = "'  -3.36%  "
$ws.Range("D11").Value = "'6.59"
$ws.Range("D12").Value = "'0.387"
$ws.Range("E12").Value = "'  +0.14%  "
$ws.Range("D13").Value = "'3.681.18"
$ws.Range("E13").Value = "'  -1.51%  "
$ws.Range("E14").Value = "'  -1.47%  "
$ws.Range("D15").Value = "'64.293.34"
$ws.Range("E15").Value = "'  -0.34%  "
$ws.Range("D16").Value = "'25.09"
$ws.Range("E16").Value = "'  -1.23%  "
$ws.Range("D17").Value = "'3.148.22"
$ws.Range("E17").Value = "'  -1.23%  "
$ws.Range("E18").Value = "'  -2.95%  "
$ws.Range("D19").Value = "'402.07"
$ws.Range("E19").Value = "'  -3.49%  "
$ws.Range("D20").Value = "'5.25"
$ws.Range("E20").Value = "'  -1.76%  "
$ws.Range("D21").Value = "'12.54"
$ws.Range("E21").Value = "'  -3.13%  "
$ws.Range("D22").Value = "'7.07"
$ws.Range("E22").Value = "'  -0.51%  "
$ws.Range("E23").Value = "'  -0.02%  "
$ws.Range("D24").Value = "'68.76"
$ws.Range("E24").Value = "'  -2.46%  "
$ws.Range("E25").Value = "'  -0.75%  "
$ws.Range("D26").Value = "'0.196"
$ws.Range("E26").Value = "'  -4.20%  "
$ws.Range("D27").Value = "'0.0000102"
$ws.Range("E27").Value = "'  -4.00%  "
$ws.Range("D28").Value = "'8.80"
$ws.Range("E28").Value = "'  +0.33%  "
$ws.Range("D29").Value = "'0.996"
$ws.Range("E29").Value = "'  -0.69%  "
$ws.Range("B30").Value = "'PancakeSwap"
$ws.Range("C30").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'1.81"
$ws.Range("E30").Value = "'  -1.58%  "
$ws.Range("B31").Value = "'USDe"
$ws.Range("C31").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "'  +0.08%  "
$ws.Range("D32").Value = "'21.27"
$ws.Range("E32").Value = "'  -2.26%  "
$ws.Range("D33").Value = "'161.20"
$ws.Range("E33").Value = "'  +1.74%  "
$ws.Range("D34").Value = "'4.87"
$ws.Range("E34").Value = "'  -4.56%  "
$ws.Range("E35").Value = "'  -1.07%  "
$ws.Range("E36").Value = "'  -2.29%  "
$ws.Range("D37").Value = "'1.35"
$ws.Range("E37").Value = "'  -1.37%  "
$ws.Range("E38").Value = "'  -1.81%  "
$ws.Range("D39").Value = "'2.638.48"
$ws.Range("E39").Value = "'  -3.54%  "
$ws.Range("E40").Value = "'  -2.95%  "
$ws.Range("E41").Value = "'  -3.21%  "
$ws.Range("D42").Value = "'38.47"
$ws.Range("E42").Value = "'  -1.85%  "
$ws.Range("E43").Value = "'  -3.86%  "
$ws.Range("D44").Value = "'0.0614"
$ws.Range("E44").Value = "'  -1.46%  "
$ws.Range("D45").Value = "'5.41"
$ws.Range("E45").Value = "'  -4.07%  "
$ws.Range("E46").Value = "'  -3.64%  "
$ws.Range("D47").Value = "'21.17"
$ws.Range("E47").Value = "'  -2.43%  "
$ws.Range("D48").Value = "'286.78"
$ws.Range("E48").Value = "'  -2.33%  "
$ws.Range("E49").Value = "'  -0.19%  "
$ws.Range("D50").Value = "'0.0977"
$ws.Range("E50").Value = "'  -1.28%  "
$ws.Range("D51").Value = "'10.48"
$ws.Range("E51").Value = "'  +0.26%  "
